$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the existing "sum" header cell (G1) onto the new
# "Save" header cell (H1) so it matches the other bold/bordered headers.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)

# Set the new header text and the corresponding value for row 2.
$ws.Range("H1").Value = "Save"
$ws.Range("H2").Value = 1
